{"js": "// Apply the edits described by the diff:\n//  - \"LUIS\" -> \"GABRIEL LIZARDI ROCHA\" (title block caption + signature line)\n//  - \"Quienes suscribimos, ... y Luis, ... Av Juan de Dios B\u00e1tiz, Juan de Dios\n//     B\u00e1tiz I, Ciudad de M\u00e9xico, CDMX, M\u00e9xico ...\" -> updated spouse name and\n//     new address\n//  - \"12 de septiembre de 2020\" -> \"16 de septiembre de 2020\"\n//  - \"Ciudad de M\u00e9xico, a 18 de junio de 2025\" -> \"Ciudad de M\u00e9xico, a 19 de\n//     junio de 2025\"\n\nconst body = context.document.body;\n\n// 1) Replace the standalone \"LUIS\" occurrences (case-sensitive, whole text of\n//    the run) \u2014 appears once in the caption block (\"Vs\" / \"LUIS\") and once in\n//    the signature block under the second signature line.\nconst luisHits = body.search(\"LUIS\", { matchCase: true, matchWholeWord: true });\nluisHits.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < luisHits.items.length; i++) {\n  luisHits.items[i].insertText(\"GABRIEL LIZARDI ROCHA\", \"Replace\");\n}\nawait context.sync();\n\n// 2) Replace the \"Quienes suscribimos...\" sentence (spouse name + address).\nconst oldIntro =\n  \"Quienes suscribimos, Zenobia Juliana Felipe Cruz y Luis, por nuestro \" +\n  \"propio derecho, se\u00f1alando como domicilio para o\u00edr y recibir \" +\n  \"notificaciones, valores y documentos, el ubicado en Av Juan de Dios \" +\n  \"B\u00e1tiz, Juan de Dios B\u00e1tiz I, Ciudad de M\u00e9xico, CDMX, M\u00e9xico, \" +\n  \"comparecemos respetuosamente para exponer:\";\nconst newIntro =\n  \"Quienes suscribimos, Zenobia Juliana Felipe Cruz y Gabriel Lizardi \" +\n  \"Rocha, por nuestro propio derecho, se\u00f1alando como domicilio para o\u00edr y \" +\n  \"recibir notificaciones, valores y documentos, el ubicado en Av 5 de \" +\n  \"Mayo 332, La Era I y II, Iztapalapa, 09720 Ciudad de M\u00e9xico, CDMX, \" +\n  \"M\u00e9xico, comparecemos respetuosamente para exponer:\";\n\nconst introHits = body.search(oldIntro, { matchCase: true });\nintroHits.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < introHits.items.length; i++) {\n  introHits.items[i].insertText(newIntro, \"Replace\");\n}\nawait context.sync();\n\n// 3) Update the marriage date.\nconst dateHits = body.search(\"12 de septiembre de 2020\", { matchCase: true });\ndateHits.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < dateHits.items.length; i++) {\n  dateHits.items[i].insertText(\"16 de septiembre de 2020\", \"Replace\");\n}\nawait context.sync();\n\n// 4) Update the closing date.\nconst closingHits = body.search(\"Ciudad de M\u00e9xico, a 18 de junio de 2025\", {\n  matchCase: true,\n});\nclosingHits.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < closingHits.items.length; i++) {\n  closingHits.items[i].insertText(\"Ciudad de M\u00e9xico, a 19 de junio de 2025\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Apply the edits described by the diff:\n#  - \"LUIS\" -> \"GABRIEL LIZARDI ROCHA\" (title block caption + signature line)\n#  - \"Quienes suscribimos, ... y Luis, ... Av Juan de Dios B\u00e1tiz, Juan de Dios\n#     B\u00e1tiz I, Ciudad de M\u00e9xico, CDMX, M\u00e9xico ...\" -> updated spouse name and\n#     new address\n#  - \"12 de septiembre de 2020\" -> \"16 de septiembre de 2020\"\n#  - \"Ciudad de M\u00e9xico, a 18 de junio de 2025\" -> \"Ciudad de M\u00e9xico, a 19 de\n#     junio de 2025\"\n\n$d = $word.ActiveDocument\n\n$wdReplaceOne = 1\n$wdReplaceAll = 2\n\n# 1) Replace the standalone \"LUIS\" occurrences (caption block + signature\n#    block) with the new spouse's name. Use MatchWholeWord so the lowercase\n#    \"Luis\" embedded in the \"Quienes suscribimos...\" sentence below is left\n#    untouched (that sentence is replaced wholesale in step 2).\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Replacement.ClearFormatting()\n$rng.Find.Execute(\"LUIS\", $true, $true, $false, $false, $false, $true, 1, $false, \"GABRIEL LIZARDI ROCHA\", $wdReplaceAll) | Out-Null\n\n# 2) Replace the \"Quienes suscribimos...\" sentence (spouse name + address).\n$oldIntro = \"Quienes suscribimos, Zenobia Juliana Felipe Cruz y Luis, por nuestro propio derecho, se\u00f1alando como domicilio para o\u00edr y recibir notificaciones, valores y documentos, el ubicado en Av Juan de Dios B\u00e1tiz, Juan de Dios B\u00e1tiz I, Ciudad de M\u00e9xico, CDMX, M\u00e9xico, comparecemos respetuosamente para exponer:\"\n$newIntro = \"Quienes suscribimos, Zenobia Juliana Felipe Cruz y Gabriel Lizardi Rocha, por nuestro propio derecho, se\u00f1alando como domicilio para o\u00edr y recibir notificaciones, valores y documentos, el ubicado en Av 5 de Mayo 332, La Era I y II, Iztapalapa, 09720 Ciudad de M\u00e9xico, CDMX, M\u00e9xico, comparecemos respetuosamente para exponer:\"\n\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Replacement.ClearFormatting()\n$rng2.Find.Execute($oldIntro, $true, $false, $false, $false, $false, $true, 1, $false, $newIntro, $wdReplaceOne) | Out-Null\n\n# 3) Update the marriage date.\n$rng3 = $d.Content\n$rng3.Find.ClearFormatting()\n$rng3.Find.Replacement.ClearFormatting()\n$rng3.Find.Execute(\"12 de septiembre de 2020\", $true, $false, $false, $false, $false, $true, 1, $false, \"16 de septiembre de 2020\", $wdReplaceOne) | Out-Null\n\n# 4) Update the closing date.\n$rng4 = $d.Content\n$rng4.Find.ClearFormatting()\n$rng4.Find.Replacement.ClearFormatting()\n$rng4.Find.Execute(\"Ciudad de M\u00e9xico, a 18 de junio de 2025\", $true, $false, $false, $false, $false, $true, 1, $false, \"Ciudad de M\u00e9xico, a 19 de junio de 2025\", $wdReplaceOne) | Out-Null\n"}
